$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '79.595.63'
$ws.Range("E2").Value = '  +3.95%  '
$ws.Range("D3").Value = '3.214.31'
$ws.Range("E3").Value = '  +6.11%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '632.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.97%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.228'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.64%  '
$ws.Range("D10").Value = '3.213.88'
$ws.Range("E10").Value = '  +6.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.585'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +35.06%  '
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.08%  '
$ws.Range("D14").Value = '3.799.18'
$ws.Range("E14").Value = '  +6.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000227'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +20.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '31.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.20%  '
$ws.Range("D17").Value = '79.299.82'
$ws.Range("E17").Value = '  +3.82%  '
$ws.Range("D18").Value = '3.194.10'
$ws.Range("E18").Value = '  +6.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.51'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +15.75%  '
$ws.Range("E22").Value = '  +25.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +17.58%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.86%  '
$ws.Range("D25").Value = '3.365.08'
$ws.Range("E25").Value = '  +7.38%  '
$ws.Range("E26").Value = '  +9.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '77.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +12.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000116'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.12%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.48'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '520.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.99'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.39%  '
$ws.Range("B36").Value = 'Cronos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.127'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +21.46%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.137'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +22.76%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.98'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +11.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.409'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '164.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '197.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.810'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +13.83%  '
$ws.Range("E47").Value = '  +8.64%  '
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '43.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.81%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.632'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.53%  '
